$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: add P1 and Q1 with header-style values 14 and 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap values in columns I, K, M, O and add P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new
    $ws.Cells.Item($r, 17).Value = 2  # Q: new
}
